# Update the "dSF" column (F) values on the active sheet to reflect
# repulled/recalculated data, per the commit message
# "repull data, push all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = -3
    3  = 2
    4  = 1
    5  = -1
    6  = -1
    7  = 1
    8  = -3
    9  = -1
    11 = 5
    13 = -3
    14 = 1
    15 = -1
    16 = 2
    17 = -1
    18 = -3
    19 = -3
    20 = 6
    21 = 0
    23 = -3
    25 = 1
    26 = -1
    27 = -2
    28 = 6
    29 = 2
    30 = 7
}

foreach ($row in $newValues.Keys) {
    $ws.Range("F$row").Value = $newValues[$row]
}
